# Auto-generated Excel COM-interop script
# Applies updated crypto Price (D) / Volume(1h) (E) figures to sheet1
# per the commit's scraped-data refresh. Source values are stored as
# plain text (prices use '.' as a thousands separator in some rows, and
# volumes are padded percent strings), so every write forces the cell's
# NumberFormat to text ("@") first -- otherwise Excel auto-coerces
# numeric-looking strings (e.g. "510.23") into Double values and mangles
# them with floating point noise (e.g. 510.23000000000002).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '57.337.19'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -2.73%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.409.83'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -3.93%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '510.23'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -4.14%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '128.57'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -4.81%  '
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.08%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.547'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -3.56%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.412.56'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -4.00%  '
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -1.50%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0951'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -6.03%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '5.17'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -4.37%  '
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -4.33%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.842.78'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -3.78%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '57.263.90'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -2.69%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '21.31'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -4.59%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.415.14'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -3.57%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '10.31'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -6.38%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '313.16'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -2.81%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.06'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -4.46%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.997'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -0.33%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.60'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -5.44%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '63.37'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -2.60%  '
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -4.62%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.999'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +0.04%  '
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -2.71%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.15'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -5.28%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '168.80'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -0.83%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0₃0715'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -6.23%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.65'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -5.28%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.14'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -5.20%  '
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +1.38%  '
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +0.01%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.999'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +0.07%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '17.63'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -4.09%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.26'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -7.68%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.83'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -4.94%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '36.28'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -2.22%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.766'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -4.15%  '
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -6.53%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '264.04'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -6.03%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '4.85'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -2.64%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.578'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -4.15%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '121.13'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -6.63%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0897'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -3.15%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0478'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -4.33%  '
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -4.22%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '16.39'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -4.83%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.682.79'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -4.45%  '

Write-Output "Updated 86 cells (Price/Volume columns) on sheet1"
